$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell for the new "time_taken" column, matching the style of the
# other header cells (bold/centered header style). Copying an existing
# header cell's formatting reuses the same cell-format index instead of
# creating a near-duplicate style.
$ws.Range("E1").Copy($ws.Range("F1"))
$ws.Range("F1").Value = "time_taken"

# Timestamp values for each data row (2-13), copied verbatim from the diff.
$timestamps = @(
    "2021-10-05 10:51:08.898633",
    "2021-10-05 10:51:08.898645",
    "2021-10-05 10:51:08.898649",
    "2021-10-05 10:51:08.898652",
    "2021-10-05 10:51:08.898656",
    "2021-10-05 10:51:08.898659",
    "2021-10-05 10:51:08.898662",
    "2021-10-05 10:51:08.898665",
    "2021-10-05 10:51:08.898668",
    "2021-10-05 10:51:08.898671",
    "2021-10-05 10:51:08.898674",
    "2021-10-05 10:51:08.898677"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
